$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "176×2=352" "813×7=5691"
Replace-Text "135×3=405" "161×7=1127"
Replace-Text "954×2=1908" "558×9=5022"
Replace-Text "713×5=3565" "595×8=4760"
Replace-Text "409×9=3681" "283×3=849"
Replace-Text "132×9=1188" "102×5=510"
Replace-Text "643×3=1929" "579×9=5211"
Replace-Text "516×4=2064" "455×5=2275"
Replace-Text "166×5=830" "819×3=2457"
Replace-Text "828×4=3312" "784×9=7056"
Replace-Text "254×6=1524" "676×5=3380"
Replace-Text "367×3=1101" "480×9=4320"
Replace-Text "254×7=1778" "206×2=412"
Replace-Text "373×9=3357" "405×5=2025"
Replace-Text "442×2=884" "549×3=1647"
Replace-Text "221×3=663" "426×5=2130"
Replace-Text "903×6=5418" "805×4=3220"
Replace-Text "326×8=2608" "568×3=1704"
Replace-Text "976×3=2928" "536×2=1072"
Replace-Text "385×3=1155" "874×4=3496"
Replace-Text "453×9=4077" "554×2=1108"
Replace-Text "231×7=1617" "344×3=1032"
Replace-Text "937×9=8433" "405×3=1215"
Replace-Text "602×6=3612" "976×7=6832"
Replace-Text "624×7=4368" "112×4=448"
